# Apply the Alvearie FHIR IG metadata refresh to the "Metadata" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version bump
$ws.Range("B3").Value = "6.0.0"

# Date refresh
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher value was previously blank
$ws.Range("B9").Value = "Alvearie Team"

# "Contact" / "No display for ContactDetail" row becomes "Jurisdiction" / "United States of America"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the now-duplicate "Contact" row entirely, shifting subsequent rows up
$ws.Rows.Item(11).Delete()

$wb.Save()
